$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3; this shifts rows 3-22 down to 4-23
$ws.Rows(3).Insert()

# Fill in the new row 3 label (matches the style used by other column-A labels)
$ws.Cells.Item(3, 1).Value = "2020-05-15 00:00:00_diff"
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(3, 1).PasteSpecial(-4122)

# Fill in the new row 3 numeric values (B..H)
$ws.Cells.Item(3, 2).Value = -2.209042350248001
$ws.Cells.Item(3, 3).Value = -1.580304150385469
$ws.Cells.Item(3, 4).Value = 3.459367396593679
$ws.Cells.Item(3, 5).Value = 4.006937868393072
$ws.Cells.Item(3, 6).Value = -0.5330921181302579
$ws.Cells.Item(3, 7).Value = -1.194266296544459
$ws.Cells.Item(3, 8).Value = 0.7471896406117955
